$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing values in rows 173-175 (columns M-Q changed)
$ws.Range("N173").Value = 12313
$ws.Range("O173").Value = 18466
$ws.Range("P173").Value = 6729
$ws.Range("Q173").Value = 9045

$ws.Range("M174").Value = 10579
$ws.Range("N174").Value = 14241
$ws.Range("O174").Value = 17900
$ws.Range("P174").Value = 7528
$ws.Range("Q174").Value = 9052

$ws.Range("M175").Value = 7259
$ws.Range("N175").Value = 13745
$ws.Range("O175").Value = 18578
$ws.Range("P175").Value = 7805
$ws.Range("Q175").Value = 8624

# Add new row 176 with the new data month (01-07-2021)
# Use a formula that evaluates to the text string, then paste-as-values,
# so Excel stores it as a plain shared string instead of auto-converting
# the date-like text into a date serial number (which would also add a
# new number-format style).
$ws.Range("A176").Formula = "=""01-07-2021"""
$ws.Range("A176").Copy()
$ws.Range("A176").PasteSpecial(-4163)
$ws.Range("B176").Value = -22646
$ws.Range("C176").Value = -15767
$ws.Range("D176").Value = 19298
$ws.Range("E176").Value = 35065
$ws.Range("F176").Value = -6879
$ws.Range("G176").Value = 4355
$ws.Range("H176").Value = 11234
$ws.Range("I176").Value = 69952
$ws.Range("J176").Value = 54363
$ws.Range("K176").Value = 185
$ws.Range("L176").Value = 763
$ws.Range("M176").Value = 5310
$ws.Range("N176").Value = 12565
$ws.Range("O176").Value = 18131
$ws.Range("P176").Value = 8964
$ws.Range("Q176").Value = 8446
$ws.Range("R176").Value = 15589
